# Add "Area" / running-area-total columns (G:H) and a compact summary
# (J:K) to the discharge worksheet, per commit "add area to Q files stn3".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 1) ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Row 2: first segment area + running totals / summary cells ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Row 3: second segment area (its own, non-shared formula) ---
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# --- Rows 4-15: remaining segment areas, entered as one shared formula ---
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Put the selection where the author left it (the new summary cells).
[void]$ws.Range("J2:K2").Select()

# Best-effort: scroll the view so column B is the left-most visible
# column (matches the saved view in the authored workbook).
$excel.ActiveWindow.ScrollColumn = 2
